$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to a number by Excel
# are forced to remain plain text (matching the inlineStr representation used
# in the workbook), then restored to the default "Normal" style so no stray
# number-format / style attribute is left behind on the cell.

$ws.Range('D2').Value = '22.462.47'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '1.570.88'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '288.28'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('E7').Value = '  +0.94%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '48.29'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.91%  '
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.133'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07484'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.937'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.24%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.887'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.37%  '
$ws.Range('D16').Value = '1.567.95'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '87.81'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.57%  '
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.352'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '16.53'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.37%  '
$ws.Range('E23').Value = '  -0.35%  '
$ws.Range('D24').Value = '22.450.84'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.387'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.569'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.46%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '153.05'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.44%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.68'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.14%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.014'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.96%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '124.28'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('D31').Value = '1.744.07'
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.049'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.793'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08312'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.47%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02463'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2269'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.88%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.06401'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.54%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.290'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.83%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.357'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.93%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.30'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.55%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.6305'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.82'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.94%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6154'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +5.50%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.773'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.056'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '125.24'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.210'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.83%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07216'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.45%  '
